$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two paragraphs involved in this edit:
#   - the paragraph ending in "...Reihenfolge)   Puru" -> becomes "Alex"
#     and gains the "_GoBack" bookmark right after the new "Alex".
#   - the paragraph ending in "...auswerten  Alex " which currently
#     carries the "_GoBack" bookmark right after "Alex"; that bookmark
#     is removed from here.
# ------------------------------------------------------------------
$targetOld = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd()
    if ($text -like "*Reihenfolge)*Puru") {
        $targetOld = $i
    }
}
if ($targetOld -eq $null) {
    throw "Could not locate target paragraph containing 'Puru'"
}

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it will be re-created at
#    the new location below). Bookmarks named "_GoBack" are hidden
#    from the normal Bookmarks collection/count, but are still
#    reachable by exact name.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ------------------------------------------------------------------
# 2) Replace "Puru" with "Alex" as the text of its own run, without
#    disturbing neighbouring runs: delete the old word, then insert
#    the new one right after what remains (this keeps the previous
#    " " run intact and creates a fresh "Alex" run in its place).
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetOld)
$pRange = $p.Range
$oldWordRange = $d.Range($pRange.End - 5, $pRange.End - 1)
if ($oldWordRange.Text -ne "Puru") {
    throw "Expected 'Puru' at computed range, found '$($oldWordRange.Text)'"
}
$oldWordRange.Delete()

$p = $d.Paragraphs.Item($targetOld)
$pRange = $p.Range
$insertPos = $pRange.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("Alex")

# ------------------------------------------------------------------
# 3) Insert the "_GoBack" bookmark immediately after the new "Alex",
#    i.e. right before the end of the paragraph. A collapsed range
#    located exactly at "end of paragraph text" is mishandled by this
#    runtime, so a throw-away character is inserted first to give the
#    bookmark a stable, non-boundary anchor point, then removed again.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetOld)
$pRange = $p.Range
$bmPos = $pRange.End - 1
$placeholder = $d.Range($bmPos, $bmPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderChar = $d.Range($bmPos, $bmPos + 1)
$placeholderChar.Delete()
